$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "73.306.17"
$ws.Range("E2").Value = "  +5.72%  "

# Row 3
$ws.Range("D3").Value = "2.664.65"
$ws.Range("E3").Value = "  +6.11%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "608.66"
$ws.Range("E5").Value = "  +2.56%  "

# Row 6
$ws.Range("D6").Value = "181.13"
$ws.Range("E6").Value = "  +2.76%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +2.80%  "

# Row 9
$ws.Range("D9").Value = "0.175"
$ws.Range("E9").Value = "  +13.76%  "

# Row 10
$ws.Range("D10").Value = "2.664.45"
$ws.Range("E10").Value = "  +6.16%  "

# Row 11
$ws.Range("E11").Value = "  +1.12%  "

# Row 12
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  +5.08%  "

# Row 13
$ws.Range("D13").Value = "5.11"
$ws.Range("E13").Value = "  +2.38%  "

# Row 14
$ws.Range("E14").Value = "  +10.01%  "

# Row 15
$ws.Range("D15").Value = "3.136.05"
$ws.Range("E15").Value = "  +5.94%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "73.157.70"
$ws.Range("E16").Value = "  +5.95%  "

# Row 17
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "27.08"
$ws.Range("E17").Value = "  +4.94%  "

# Row 18
$ws.Range("D18").Value = "2.661.53"
$ws.Range("E18").Value = "  +6.21%  "

# Row 19
$ws.Range("D19").Value = "387.12"
$ws.Range("E19").Value = "  +7.14%  "

# Row 20
$ws.Range("D20").Value = "11.65"
$ws.Range("E20").Value = "  +6.33%  "

# Row 21
$ws.Range("D21").Value = "7.98"
$ws.Range("E21").Value = "  +5.43%  "

# Row 22
$ws.Range("E22").Value = "  +4.26%  "

# Row 23
$ws.Range("E23").Value = "  +22.68%  "

# Row 24
$ws.Range("D24").Value = "73.85"
$ws.Range("E24").Value = "  +5.37%  "

# Row 25
$ws.Range("D25").Value = "4.47"
$ws.Range("E25").Value = "  +6.63%  "

# Row 26
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").Value = "9.98"
$ws.Range("E27").Value = "  +10.93%  "

# Row 28
$ws.Range("D28").Value = "2.802.57"
$ws.Range("E28").Value = "  +6.57%  "

# Row 29
$ws.Range("E29").Value = "  -0.28%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0981"
$ws.Range("E30").Value = "  +10.98%  "

# Row 31
$ws.Range("D31").Value = "540.01"
$ws.Range("E31").Value = "  +5.60%  "

# Row 32
$ws.Range("D32").Value = "8.13"
$ws.Range("E32").Value = "  +5.38%  "

# Row 33
$ws.Range("E33").Value = "  +10.36%  "

# Row 34
$ws.Range("E34").Value = "  +4.65%  "

# Row 35
$ws.Range("E35").Value = "  +0.03%  "

# Row 36
$ws.Range("D36").Value = "164.13"
$ws.Range("E36").Value = "  +1.47%  "

# Row 37
$ws.Range("D37").Value = "19.45"
$ws.Range("E37").Value = "  +4.34%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.114"
$ws.Range("E38").Value = "  -4.07%  "

# Row 39
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "1.42"
$ws.Range("E39").Value = "  +9.51%  "

# Row 40
$ws.Range("E40").Value = "  +2.32%  "

# Row 41
$ws.Range("D41").Value = "1.86"
$ws.Range("E41").Value = "  +9.35%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "5.17"
$ws.Range("E42").Value = "  +8.20%  "

# Row 43
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.67"
$ws.Range("E43").Value = "  +15.99%  "

# Row 44
$ws.Range("E44").Value = "  +0.10%  "

# Row 45
$ws.Range("D45").Value = "0.338"
$ws.Range("E45").Value = "  +5.80%  "

# Row 46
$ws.Range("D46").Value = "39.84"
$ws.Range("E46").Value = "  +2.95%  "

# Row 47
$ws.Range("D47").Value = "152.56"
$ws.Range("E47").Value = "  +1.92%  "

# Row 48
$ws.Range("D48").Value = "3.73"
$ws.Range("E48").Value = "  +5.00%  "

# Row 49
$ws.Range("E49").Value = "  +6.77%  "

# Row 50
$ws.Range("E50").Value = "  +10.83%  "

# Row 51
$ws.Range("D51").Value = "0.0₆0270"
$ws.Range("E51").Value = "  +10.46%  "
